# Update scripts with new TPM values.
#
# The underlying dataset (Adam2 -> Itga6 signalling) was recomputed with an
# updated TPM table. This changes the previously-computed statistics for the
# existing "FAPs" sending-cluster rows (row 2-6), and adds a mirrored block of
# rows (7-11) for a new "MuSCs" sending cluster that sends the same
# ligand/receptor pair to the same five target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D are strings (cluster / symbol names), E-T are numeric stats.
$stringCols = @("A", "B", "C", "D")
$numericCols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# One entry per data row, in sheet order (row 2 .. row 11).
$rows = @(
    @{ Row = 2;  A = "FAPs";  B = "Adam2"; C = "Itga6"; D = "ECs";              E = 1; F = 0.3333333333333333; G = 0.1251886666666667; H = 0.375566; I = 0.6155719715657366; J = 0.7060466830097307; K = 3; L = 1; M = 145.7087706666667;  N = 437.126312;   O = 0.5445232453600627;  P = 0.5461141113270247;  Q = 18.24108672139911;   R = 164.169780492592;  S = 0.3351932477096671;  T = 0.3855820568472526 },
    @{ Row = 3;  A = "FAPs";  B = "Adam2"; C = "Itga6"; D = "FAPs";             E = 1; F = 0.3333333333333333; G = 0.1251886666666667; H = 0.375566; I = 0.6155719715657366; J = 0.7060466830097307; K = 3; L = 1; M = 0.896351;             N = 2.689053;     O = 0.003349722554576428; P = 0.003359509023117945; Q = 0.1122129865553333;  R = 1.009916878998;    S = 0.002061995317118828; T = 0.002371970202313685 },
    @{ Row = 4;  A = "FAPs";  B = "Adam2"; C = "Itga6"; D = "Inflammatory-Mac"; E = 1; F = 0.3333333333333333; G = 0.1251886666666667; H = 0.375566; I = 0.6155719715657366; J = 0.7060466830097307; K = 3; L = 1; M = 51.59199533333333;   N = 154.775986;   O = 0.1928026748491032;  P = 0.1933659624890163;  Q = 6.458733106452889;  R = 58.128597958076;   S = 0.1186839226800101;  T = 0.136525396422354  },
    @{ Row = 5;  A = "FAPs";  B = "Adam2"; C = "Itga6"; D = "MuSCs";            E = 1; F = 0.3333333333333333; G = 0.1251886666666667; H = 0.375566; I = 0.6155719715657366; J = 0.7060466830097307; K = 2; L = 1; M = 2.338518;             N = 4.677036;     O = 0.0087391953474509;  P = 0.005843151713055659; Q = 0.292755950396;      R = 1.756535702376;    S = 0.005379603709928463; T = 0.004125537885325574 },
    @{ Row = 6;  A = "FAPs";  B = "Adam2"; C = "Itga6"; D = "Resolving-Mac";    E = 1; F = 0.3333333333333333; G = 0.1251886666666667; H = 0.375566; I = 0.6155719715657366; J = 0.7060466830097307; K = 3; L = 1; M = 67.05398933333333;   N = 201.161968;   O = 0.2505851618888069;  P = 0.2513172654477853;  Q = 8.394399519320888;  R = 75.54959567388801; S = 0.1542532021490121;  T = 0.1774417216524848 },
    @{ Row = 7;  A = "MuSCs"; B = "Adam2"; C = "Itga6"; D = "ECs";              E = 1; F = 0.5;                 G = 0.078181;           H = 0.156362; I = 0.3844280284342634; J = 0.2939533169902694; K = 3; L = 1; M = 145.7087706666667;  N = 437.126312;   O = 0.5445232453600627;  P = 0.5461141113270247;  Q = 11.39165739949066;  R = 68.349944396944;   S = 0.2093299976503956;  T = 0.1605320544797721 },
    @{ Row = 8;  A = "MuSCs"; B = "Adam2"; C = "Itga6"; D = "FAPs";             E = 1; F = 0.5;                 G = 0.078181;           H = 0.156362; I = 0.3844280284342634; J = 0.2939533169902694; K = 3; L = 1; M = 0.896351;             N = 2.689053;     O = 0.003349722554576428; P = 0.003359509023117945; Q = 0.070077617531;     R = 0.420465705186;    S = 0.001287727237457601; T = 0.0009875388208042595 },
    @{ Row = 9;  A = "MuSCs"; B = "Adam2"; C = "Itga6"; D = "Inflammatory-Mac"; E = 1; F = 0.5;                 G = 0.078181;           H = 0.156362; I = 0.3844280284342634; J = 0.2939533169902694; K = 3; L = 1; M = 51.59199533333333;   N = 154.775986;   O = 0.1928026748491032;  P = 0.1933659624890163;  Q = 4.033513787155333;  R = 24.201082722932;   S = 0.07411875216909307; T = 0.05684056606666234 },
    @{ Row = 10; A = "MuSCs"; B = "Adam2"; C = "Itga6"; D = "MuSCs";            E = 1; F = 0.5;                 G = 0.078181;           H = 0.156362; I = 0.3844280284342634; J = 0.2939533169902694; K = 2; L = 1; M = 2.338518;             N = 4.677036;     O = 0.0087391953474509;  P = 0.005843151713055659; Q = 0.182827675758;      R = 0.7313107030320001; S = 0.003359591637522437; T = 0.001717613827730086 },
    @{ Row = 11; A = "MuSCs"; B = "Adam2"; C = "Itga6"; D = "Resolving-Mac";    E = 1; F = 0.5;                 G = 0.078181;           H = 0.156362; I = 0.3844280284342634; J = 0.2939533169902694; K = 3; L = 1; M = 67.05398933333333;   N = 201.161968;   O = 0.2505851618888069;  P = 0.2513172654477853;  Q = 5.242347940069333;  R = 31.454087640416;   S = 0.09633195973979473; T = 0.07387554379530049 }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    foreach ($col in $stringCols) {
        $ws.Range("$col$r").Value = $entry[$col]
    }

    foreach ($col in $numericCols) {
        $ws.Range("$col$r").Value = $entry[$col]
    }
}
